$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.1111111111111111
$ws.Range("C2").Value2 = 0.7777777777777778
$ws.Range("P2").Value2 = 0.1111111111111111
$ws.Range("P3").Value2 = 1
$ws.Range("P4").Value2 = 1
$ws.Range("F6").Value2 = 0.125
$ws.Range("J6").Value2 = 0.375
$ws.Range("R6").Value2 = 0.25
$ws.Range("S6").Value2 = 0.25
$ws.Range("D7").Value2 = 0.1428571428571428
$ws.Range("F7").Value2 = 0.1428571428571428
$ws.Range("Q7").Value2 = 0.4285714285714285
$ws.Range("S7").Value2 = 0.2857142857142857
$ws.Range("B8").Value2 = 0.2
$ws.Range("O8").Value2 = 0.2
$ws.Range("Q8").Value2 = 0.1
$ws.Range("R8").Value2 = 0.3
$ws.Range("S8").Value2 = 0.2
$ws.Range("B9").Value2 = 0.1666666666666667
$ws.Range("J9").Value2 = 0.1666666666666667
$ws.Range("O9").Value2 = 0.1666666666666667
$ws.Range("S9").Value2 = 0.5
$ws.Range("B10").Value2 = 0.09803921568627451
$ws.Range("D10").Value2 = 0.0196078431372549
$ws.Range("F10").Value2 = 0.05882352941176471
$ws.Range("J10").Value2 = 0.07843137254901961
$ws.Range("Q10").Value2 = 0.3529411764705883
$ws.Range("R10").Value2 = 0.1372549019607843
$ws.Range("S10").Value2 = 0.2549019607843137
$ws.Range("G11").Value2 = 0.5
$ws.Range("K11").Value2 = 0.5
$ws.Range("G13").Value2 = 1
$ws.Range("F15").Value2 = 0.1
$ws.Range("H15").Value2 = 0.2
$ws.Range("J15").Value2 = 0.4
$ws.Range("S15").Value2 = 0.3
$ws.Range("J16").Value2 = 0.6666666666666666
$ws.Range("O16").Value2 = 0.1111111111111111
$ws.Range("S16").Value2 = 0.2222222222222222
$ws.Range("F17").Value2 = 0.04761904761904762
$ws.Range("H17").Value2 = 0.04761904761904762
$ws.Range("I17").Value2 = 0.1428571428571428
$ws.Range("J17").Value2 = 0.5714285714285714
$ws.Range("K17").Value2 = 0.04761904761904762
$ws.Range("M17").Value2 = 0.04761904761904762
$ws.Range("O17").Value2 = 0.04761904761904762
$ws.Range("S17").Value2 = 0.04761904761904762
$ws.Range("H18").Value2 = 0.1666666666666667
$ws.Range("I18").Value2 = 0.08333333333333333
$ws.Range("J18").Value2 = 0.5
$ws.Range("M18").Value2 = 0.08333333333333333
$ws.Range("S18").Value2 = 0.1666666666666667
$ws.Range("H19").Value2 = 0.1612903225806452
$ws.Range("I19").Value2 = 0.06451612903225806
$ws.Range("J19").Value2 = 0.5161290322580645
$ws.Range("K19").Value2 = 0.06451612903225806
$ws.Range("M19").Value2 = 0.03225806451612903
$ws.Range("O19").Value2 = 0.1290322580645161
$ws.Range("S19").Value2 = 0.03225806451612903
